$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing (pre-edit) software/output locations so we can
# re-insert them further down the sheet before overwriting B2/B3 with the
# new JTMT-based paths.
$oldSoftwareLocation = $ws.Range("B2").Text
$oldOutputLocation   = $ws.Range("B3").Text

# New paths (local JTMT working copy) replace the old network-share paths.
$ws.Range("B2").Value = "C:\Users\dpere\Documents\JTMT\forecast\create_forecast_basic\current"
$ws.Range("B3").Value = "C:\Users\dpere\Documents\JTMT\forecast_by_version\V4\BASE_YEAR"

# Leave a couple of blank spacer rows, matching the target layout.
$ws.Range("B5").NumberFormat = "General"
$ws.Range("B6").NumberFormat = "General"

# Re-add the old paths further down the sheet for reference.
$ws.Range("B8").Value = $oldSoftwareLocation
$ws.Range("B9").Value = $oldOutputLocation

# Update the selected range to match the new edit location.
$ws.Range("B2:B3").Select() | Out-Null
